$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Modelo de dados/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$pRange = $para.Range
Write-Output "pRange start=$($pRange.Start) end=$($pRange.End) text=[$($pRange.Text)]"

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="InfoBlue"/>
<w:jc w:val="center"/>
<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:i w:val="0"/><w:color w:val="00000A"/><w:lang w:eastAsia="en-US"/></w:rPr>
</w:pPr>
<w:r>
<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:i w:val="0"/><w:color w:val="00000A"/><w:lang w:eastAsia="en-US"/></w:rPr>
<w:t>https://svn.mec.gov.br/simec/simec/trunk/docs/01-Especificacao/Banco de dados/Modelo de dados/</w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:i w:val="0"/><w:color w:val="00000A"/><w:lang w:eastAsia="en-US"/></w:rPr>
<w:t>mer</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:i w:val="0"/><w:color w:val="00000A"/><w:lang w:eastAsia="en-US"/></w:rPr>
<w:t>_/</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$pRange.InsertXML($xmlFrag)
Write-Output "Done"
